# Swap the order of the first two slides ("Rotation" title slide and
# "첫 번째 튜토리얼" slide) in the presentation.
$p = $ppt.ActivePresentation

# Move what is currently the 2nd slide ("첫 번째 튜토리얼") so it becomes
# the 1st slide; the former 1st slide ("Rotation") shifts down to 2nd.
$s2 = $p.Slides.Item(2)
$s2.MoveTo(1)
